$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header in H2 and value in H3
$ws.Range("H2").Value = "Words"
$ws.Range("H3").Value = 8370

# Trim updates in column F
$ws.Range("F13").Value = 2.95
$ws.Range("F26").Value = 2.1
$ws.Range("F28").Value = 2.9

# Update selection to F29
$ws.Range("F29").Select()
